# Update town close year columns (#769)
# Rename "2023/2024" year-specific headers to generic "Prior Year / Curr. Year"
# headers on Sheet1's header row, and widen the affected columns to fit the
# new (longer) labels.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Header text updates -----------------------------------------------
# Cells are written in this specific order so that newly-introduced shared
# strings land in the same relative order as the target workbook.
$ws.Range("J1").Value  = "Prior Year Oc %"
$ws.Range("L1").Value  = "Prior Year % of Curr. Year"
$ws.Range("K1").Value  = "Curr. Year Occ %"
$ws.Range("N1").Value  = "Prior Year DWELDAT VAL"
$ws.Range("P1").Value  = "Prior Year LMV"
$ws.Range("Q1").Value  = "Prior Year BMV"
$ws.Range("R1").Value  = "Prior Year Total MV"
$ws.Range("S1").Value  = "Prior Year LAV"
$ws.Range("T1").Value  = "Prior Year BAV"
$ws.Range("U1").Value  = "Prior Year Total AV"
$ws.Range("V1").Value  = "Curr. Year LMV"
$ws.Range("W1").Value  = "Curr. Year BMV"
$ws.Range("X1").Value  = "Curr. Year Total MV"
$ws.Range("Y1").Value  = "Curr. Year LAV"
$ws.Range("Z1").Value  = "Curr. Year BAV"
$ws.Range("AA1").Value = "Curr. Year Total AV"
$ws.Range("O1").Value  = "Curr. Year DWELDAT VAL"

# --- Column width updates -----------------------------------------------
# Widen the columns whose headers grew so the new labels fit (matches the
# "Fix column widths" / "Extend column width" / "Declutter column widths"
# commits). Only the columns whose header text actually changed are
# touched here.
$ws.Columns.Item(10).ColumnWidth = 16.053385416666668   # J  Prior Year Oc %
$ws.Columns.Item(11).ColumnWidth = 16.608072916666668   # K  Curr. Year Occ %
$ws.Columns.Item(12).ColumnWidth = 24.166666666666668   # L  Prior Year % of Curr. Year
$ws.Columns.Item(14).ColumnWidth = 24.498697916666668   # N  Prior Year DWELDAT VAL
$ws.Columns.Item(15).ColumnWidth = 24.830729166666668   # O  Curr. Year DWELDAT VAL
$ws.Columns.Item(16).ColumnWidth = 15.944010416666666   # P  Prior Year LMV
$ws.Columns.Item(17).ColumnWidth = 16.608072916666668   # Q  Prior Year BMV
$ws.Columns.Item(18).ColumnWidth = 20.166666666666668   # R  Prior Year Total MV
$ws.Columns.Item(19).ColumnWidth = 16.053385416666668   # S  Prior Year LAV
$ws.Columns.Item(20).ColumnWidth = 15.498697916666666   # T  Prior Year BAV
$ws.Columns.Item(21).ColumnWidth = 19.276041666666668   # U  Prior Year Total AV
$ws.Columns.Item(22).ColumnWidth = 16.498697916666668   # V  Curr. Year LMV
$ws.Columns.Item(23).ColumnWidth = 15.944010416666666   # W  Curr. Year BMV
$ws.Columns.Item(24).ColumnWidth = 20.830729166666668   # X  Curr. Year Total MV
$ws.Columns.Item(25).ColumnWidth = 15.721354166666666   # Y  Curr. Year LAV
$ws.Columns.Item(26).ColumnWidth = 15.276041666666666   # Z  Curr. Year BAV
$ws.Columns.Item(27).ColumnWidth = 17.385416666666668   # AA Curr. Year Total AV
